$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells whose new values look numeric,
# so Excel keeps them as text (matching original formatting like "243.00", "0.620", etc.)
$textCells = @(
    'D5', 'D6', 'D8', 'D9', 'D10', 'D11', 'D14', 'D16', 'D17', 'D18', 'D19', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D30', 'D31', 'D32', 'D33', 'D34', 'D37', 'D39', 'D40', 'D41', 'D42', 'D44', 'D45', 'D48', 'D50'
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '33.986.19'
$ws.Range('E2').Value = '  -1.98%  '
$ws.Range('D3').Value = '1.779.59'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').Value = '220.33'
$ws.Range('E5').Value = '  -2.37%  '
$ws.Range('D6').Value = '0.546'
$ws.Range('E6').Value = '  -1.65%  '
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').Value = '31.27'
$ws.Range('E8').Value = '  -4.63%  '
$ws.Range('D9').Value = '0.285'
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('D10').Value = '0.0701'
$ws.Range('E10').Value = '  +3.63%  '
$ws.Range('D11').Value = '0.0922'
$ws.Range('E11').Value = '  -1.57%  '
$ws.Range('D12').Value = '2.040.46'
$ws.Range('E12').Value = '  -0.77%  '
$ws.Range('D13').Value = '1.779.15'
$ws.Range('E13').Value = '  -1.08%  '
$ws.Range('D14').Value = '10.56'
$ws.Range('E14').Value = '  -5.39%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '34.012.37'
$ws.Range('E15').Value = '  -1.73%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = '0.620'
$ws.Range('E16').Value = '  -3.09%  '
$ws.Range('D17').Value = '4.19'
$ws.Range('E17').Value = '  -2.35%  '
$ws.Range('D18').Value = '67.76'
$ws.Range('E18').Value = '  -2.77%  '
$ws.Range('D19').Value = '243.00'
$ws.Range('E19').Value = '  -5.19%  '
$ws.Range('D20').Value = '0.0₃0768'
$ws.Range('E20').Value = '  -1.43%  '
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('D22').Value = '10.53'
$ws.Range('E22').Value = '  +0.57%  '
$ws.Range('D23').Value = '4.06'
$ws.Range('E23').Value = '  -4.35%  '
$ws.Range('D24').Value = '2.11'
$ws.Range('E24').Value = '  -0.78%  '
$ws.Range('D25').Value = '157.72'
$ws.Range('E25').Value = '  -0.39%  '
$ws.Range('D26').Value = '16.24'
$ws.Range('E26').Value = '  -1.75%  '
$ws.Range('D27').Value = '6.96'
$ws.Range('E27').Value = '  -2.65%  '
$ws.Range('D28').Value = '0.111'
$ws.Range('E28').Value = '  -2.62%  '
$ws.Range('E29').Value = '  +0.27%  '
$ws.Range('D30').Value = '0.0516'
$ws.Range('E30').Value = '  -0.70%  '
$ws.Range('D31').Value = '1.19'
$ws.Range('E31').Value = '  -0.48%  '
$ws.Range('D32').Value = '3.66'
$ws.Range('E32').Value = '  -3.82%  '
$ws.Range('D33').Value = '3.49'
$ws.Range('E33').Value = '  -3.24%  '
$ws.Range('D34').Value = '1.82'
$ws.Range('E34').Value = '  -4.45%  '
$ws.Range('D35').Value = '1.394.06'
$ws.Range('E35').Value = '  -4.22%  '
$ws.Range('E36').Value = '  -1.09%  '
$ws.Range('D37').Value = '0.622'
$ws.Range('E37').Value = '  -2.50%  '
$ws.Range('E38').Value = '  -3.09%  '
$ws.Range('D39').Value = '0.923'
$ws.Range('E39').Value = '  +2.04%  '
$ws.Range('B40').Value = 'HuobiToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D40').Value = '2.34'
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '78.78'
$ws.Range('E41').Value = '  -5.26%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Value = '2.70'
$ws.Range('E42').Value = '  -5.34%  '
$ws.Range('E43').Value = '  +0.42%  '
$ws.Range('D44').Value = '0.0490'
$ws.Range('E44').Value = '  -3.72%  '
$ws.Range('D45').Value = '5.82'
$ws.Range('E45').Value = '  -1.96%  '
$ws.Range('E46').Value = '  -1.11%  '
$ws.Range('D47').Value = '1.932.73'
$ws.Range('E47').Value = '  -1.20%  '
$ws.Range('D48').Value = '105.03'
$ws.Range('E48').Value = '  +3.48%  '
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('D50').Value = '11.70'
$ws.Range('E50').Value = '  -2.70%  '
$ws.Range('E51').Value = '  +0.43%  '
